$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.21729837854633161
$ws.Cells.Item(2, 1).Value = -0.0059999999616131561
$ws.Cells.Item(3, 1).Value = -0.0039999999661350927
$ws.Cells.Item(4, 1).Value = -0.0079999999378230768
$ws.Cells.Item(5, 1).Value = -0.0029999999656800114
$ws.Cells.Item(6, 1).Value = -0.0019999999618800501
$ws.Cells.Item(7, 1).Value = -0.009999999912023938
$ws.Cells.Item(8, 1).Value = -0.0099999999123117078
$ws.Cells.Item(9, 1).Value = 0.055239506568542041
$ws.Cells.Item(10, 1).Value = -0.0019999999581035155
$ws.Cells.Item(11, 1).Value = -0.0029999999515375464
$ws.Cells.Item(12, 1).Value = -0.0034999999476328369
$ws.Cells.Item(13, 1).Value = -0.012786368774028389
$ws.Cells.Item(14, 1).Value = -0.0079999999130215826
$ws.Cells.Item(15, 1).Value = -0.00099999995581523393
$ws.Cells.Item(16, 1).Value = -0.0019999999490094567
$ws.Cells.Item(17, 1).Value = -0.0019999999480209141
$ws.Cells.Item(18, 1).Value = -0.0039999999353890203
$ws.Cells.Item(19, 1).Value = -0.0039999999728750346
$ws.Cells.Item(20, 1).Value = -0.059533955893755319
$ws.Cells.Item(21, 1).Value = -0.0039999999581210588
$ws.Cells.Item(22, 1).Value = -0.0039999999578590462
$ws.Cells.Item(23, 1).Value = -0.0049999999581844534
$ws.Cells.Item(24, 1).Value = -0.01999999986113199
$ws.Cells.Item(25, 1).Value = -0.019999999859209971
$ws.Cells.Item(26, 1).Value = -0.0024999999572870024
$ws.Cells.Item(27, 1).Value = -0.002499999955343668
$ws.Cells.Item(28, 1).Value = -0.0019999999485360576
$ws.Cells.Item(29, 1).Value = -0.0069999999108141253
$ws.Cells.Item(30, 1).Value = -0.059999999578964847
$ws.Cells.Item(31, 1).Value = -0.0069999999078049768
$ws.Cells.Item(32, 1).Value = 0.016260856065867912
$ws.Cells.Item(33, 1).Value = -0.0039999999265916131
